$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.029.98"
$ws.Range("E2").Value = "  +5.33%  "
$ws.Range("D3").Value = "2.317.39"
$ws.Range("E3").Value = "  +2.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Formula = "'518.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.68%  "
$ws.Range("D6").Formula = "'134.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.45%  "
$ws.Range("D7").Formula = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Formula = "'0.538"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").Value = "2.340.65"
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("E10").Value = "  +9.23%  "
$ws.Range("D11").Formula = "'0.155"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("E12").Value = "  +6.32%  "
$ws.Range("D13").Formula = "'0.345"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.12%  "
$ws.Range("D14").Formula = "'24.04"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.04%  "
$ws.Range("D15").Value = "2.731.09"
$ws.Range("E15").Value = "  +2.57%  "
$ws.Range("D16").Value = "56.864.88"
$ws.Range("E17").Value = "  +5.07%  "
$ws.Range("D18").Value = "2.333.97"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("E19").Value = "  +3.31%  "
$ws.Range("D20").Formula = "'4.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.80%  "
$ws.Range("D21").Formula = "'321.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.24%  "
$ws.Range("E22").Value = "  +4.99%  "
$ws.Range("D23").Formula = "'0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").Formula = "'61.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Formula = "'0.994"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Formula = "'0.159"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.90%  "
$ws.Range("D27").Formula = "'7.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.97%  "
$ws.Range("D28").Formula = "'172.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  +12.53%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +7.05%  "
$ws.Range("E31").Value = "  +5.34%  "
$ws.Range("E32").Value = "  +4.89%  "
$ws.Range("D33").Formula = "'18.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Formula = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Formula = "'0.951"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("E37").Value = "  +5.83%  "
$ws.Range("E38").Value = "  +9.06%  "
$ws.Range("E39").Value = "  +8.99%  "
$ws.Range("D40").Formula = "'37.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.71%  "
$ws.Range("D41").Formula = "'0.383"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.22%  "
$ws.Range("D42").Formula = "'140.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +13.22%  "
$ws.Range("E43").Value = "  +7.47%  "
$ws.Range("D44").Formula = "'279.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +15.99%  "
$ws.Range("D45").Formula = "'5.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.62%  "
$ws.Range("D46").Formula = "'0.0511"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.03%  "
$ws.Range("D47").Formula = "'0.0931"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.27%  "
$ws.Range("D48").Formula = "'0.563"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("E50").Value = "  +6.10%  "
$ws.Range("D51").Formula = "'16.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.11%  "
